# "aggiornamento a 9/09 compreso" -- extend the daily series through
# 2021-09-09 (date serial 44448), adding rows 367-374 after the existing
# last row (366, serial 44440 / 2021-09-01).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 366
$firstNew = 367
$lastNew = 374

# Clone column A's formatting (date number format / style) from the last
# existing row down into the new rows, then fill in the real values.
$ws.Range("A$srcRow").Copy()
$ws.Range("A$($firstNew):A$($lastNew)").PasteSpecial(-4122)

for ($r = $firstNew; $r -le $lastNew; $r++) {
    $serial = 44440 + ($r - $srcRow)
    $ws.Cells.Item($r, 1).Value = $serial
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}
